# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 26 (pushing the existing
# rows 26-41 down to 27-42) and populate it with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("26:26").Insert()

$ws.Range("A26").Value = 10
$ws.Range("B26").Value = "Vega Modelo de Temuco"
$ws.Range("C26").Value = "La Araucanía"
$ws.Range("D26").Value = 44726
$ws.Range("E26").Value = 9
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100108
$ws.Range("H26").Value = "Tropicales y subtropicales"
$ws.Range("I26").Value = 100108003
$ws.Range("J26").Value = "Maracuyá"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 30
$ws.Range("N26").Value = 34000
$ws.Range("O26").Value = 34000
$ws.Range("P26").Value = 34000
$ws.Range("Q26").Value = "$/caja 18 kilos"
$ws.Range("R26").Value = "Región de Arica y Parinacota"
$ws.Range("S26").Value = 1889
$ws.Range("T26").Value = 18
